# poprawa wyswietlania list i dodanie watkow dla walut i surowcow
#
# 1. Remove the first paragraph entirely ("Przy kazdym zakupie...").
# 2. Split the "Monitory w klasie MenuFunctionality" paragraph into three:
#       "Dodanie indeksow"
#       "Poprawa tabel w okienkach z informacjami"
#       "Wywalenie zbednych metod z interfejsu " + "allinstances"
#    (the spell-checker wrap around the former "MenuFunctionality" run is
#    reused/renamed in place so the proofErr markers survive intact).
# 3. Replace the "Przeniesienie buy z inwestora do menu functionality"
#    paragraph's text with "Watek surowcow i walut", keeping the bookmark.

$d = $word.ActiveDocument

# --- 1. Drop the leading paragraph -----------------------------------------
$d.Paragraphs.Item(1).Range.Delete()

# --- 2. Split "Monitory w klasie MenuFunctionality" ------------------------
# Break the paragraph right after the "Monitory w klasie " prefix, so the
# spell-check-wrapped "MenuFunctionality" run ends up alone in its own
# paragraph (its proofErr wrap is preserved by the split).
$splitRng = $d.Content.Duplicate
$splitRng.Find.Execute("Monitory w klasie ")
$splitRng.InsertParagraphAfter()

# Rename the now-isolated prefix paragraph's text.
$d.Content.Find.Execute("Monitory w klasie ", $false, $false, $false, $false, `
    $false, $true, 1, $false, "Dodanie indeksów", 2)

# Insert the brand new "Poprawa tabel..." paragraph right after it.
$prefixPara = $d.Paragraphs.Item(2)
$prefixPara.Range.InsertParagraphAfter()
$d.Paragraphs.Item(3).Range.Text = "Poprawa tabel w okienkach z informacjami"

# Prepend the new lead-in text to the paragraph still holding the
# (proofErr-wrapped) former "MenuFunctionality" run.
$menuPara = $d.Paragraphs.Item(4)
$menuPara.Range.InsertBefore("Wywalenie zbędnych metod z interfejsu ")

# Rename that wrapped word in place -- the spellStart/spellEnd proofErr
# markers stay exactly where they were, just around the new word.
$d.Content.Find.Execute("MenuFunctionality", $false, $false, $false, $false, `
    $false, $true, 1, $false, "allinstances", 2)

# --- 3. Replace the "Przeniesienie ... menu functionality" paragraph -------
$d.Content.Find.Execute("Przeniesienie buy z inwestora do menu functionality", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "Wątek surowców i walut", 2)
